$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 355.4
$ws.Range("I28").Value = 253.75
$ws.Range("K28").Value = 253.75
$ws.Range("M28").Value = 231.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 196.8421
$ws.Range("I33").Value = 222.5
$ws.Range("K33").Value = 222.5
$ws.Range("M33").Value = 6.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2977224.8
$ws.Range("I112").Value = 62500400
$ws.Range("J112").Value = 1066.025
$ws.Range("K112").Value = 187501200
$ws.Range("L112").Value = 3198.075
$ws.Range("M112").Value = -187500092
$ws.Range("N112").Value = -5414.075000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 11908667
$ws.Range("I116").Value = 22728628
$ws.Range("K116").Value = 22728628
$ws.Range("M116").Value = -22725186

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 743.4651
$ws.Range("I129").Value = 296.8
$ws.Range("J129").Value = 802.2368
$ws.Range("K129").Value = 890.4000000000001
$ws.Range("L129").Value = 2406.7104
$ws.Range("M129").Value = 4109.6
$ws.Range("N129").Value = -12406.7104

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 49780
$ws.Range("J133").Value = 49780
$ws.Range("L133").Value = 49780
$ws.Range("N133").Value = -59900

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2099.7837
$ws.Range("I138").Value = 1109.0333
$ws.Range("J138").Value = 2775.2954
$ws.Range("K138").Value = 3327.0999
$ws.Range("L138").Value = 8325.886200000001
$ws.Range("M138").Value = 1812.9001
$ws.Range("N138").Value = -18605.8862

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2547.125
$ws.Range("I141").Value = 2170.75
$ws.Range("J141").Value = 3676.25
$ws.Range("K141").Value = 6512.25
$ws.Range("L141").Value = 11028.75
$ws.Range("M141").Value = -1332.25
$ws.Range("N141").Value = -21388.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5102.47
$ws.Range("I32").Value = 4815.8047
$ws.Range("J32").Value = 6408.3887
$ws.Range("K32").Value = 4815.8047
$ws.Range("L32").Value = 6408.3887
$ws.Range("M32").Value = -4528.8047
$ws.Range("N32").Value = -6982.3887

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2570.875
$ws.Range("I45").Value = 2535.2144
$ws.Range("J45").Value = 2620.8
$ws.Range("K45").Value = 2535.2144
$ws.Range("L45").Value = 2620.8
$ws.Range("M45").Value = -2158.2144
$ws.Range("N45").Value = -3374.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1811.075
$ws.Range("I61").Value = 1516.7188
$ws.Range("J61").Value = 2988.5
$ws.Range("K61").Value = 1516.7188
$ws.Range("L61").Value = 2988.5
$ws.Range("M61").Value = -1304.7188
$ws.Range("N61").Value = -3412.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1811.075
$ws.Range("I136").Value = 1516.7188
$ws.Range("J136").Value = 2988.5
$ws.Range("K136").Value = 4550.1564
$ws.Range("L136").Value = 8965.5
$ws.Range("M136").Value = -2000.1564
$ws.Range("N136").Value = -14065.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1002
$ws.Range("I8").Value = 1002
$ws.Range("K8").Value = 1002
$ws.Range("M8").Value = -862

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3039.1592
$ws.Range("I134").Value = 3121.6428
$ws.Range("K134").Value = 9364.928400000001
$ws.Range("M134").Value = -6829.928400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3482.678
$ws.Range("I31").Value = 1670.7241
$ws.Range("J31").Value = 5234.2334
$ws.Range("K31").Value = 1670.7241
$ws.Range("L31").Value = 5234.2334
$ws.Range("M31").Value = -1375.7241
$ws.Range("N31").Value = -5824.2334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3482.678
$ws.Range("I34").Value = 1670.7241
$ws.Range("J34").Value = 5234.2334
$ws.Range("K34").Value = 1670.7241
$ws.Range("L34").Value = 5234.2334
$ws.Range("M34").Value = -1468.7241
$ws.Range("N34").Value = -5638.2334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 17796.129
$ws.Range("I58").Value = 1958.75
$ws.Range("J58").Value = 27798.684
$ws.Range("K58").Value = 1958.75
$ws.Range("L58").Value = 27798.684
$ws.Range("M58").Value = -1755.75
$ws.Range("N58").Value = -28204.684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3910.7827
$ws.Range("I99").Value = 2767.5293
$ws.Range("J99").Value = 7150
$ws.Range("K99").Value = 2767.5293
$ws.Range("L99").Value = 7150
$ws.Range("M99").Value = -1269.5293
$ws.Range("N99").Value = -10146

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1107.75
$ws.Range("I107").Value = 441.73334
$ws.Range("K107").Value = 441.73334
$ws.Range("M107").Value = 1478.26666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3910.7827
$ws.Range("I126").Value = 2767.5293
$ws.Range("J126").Value = 7150
$ws.Range("K126").Value = 8302.5879
$ws.Range("L126").Value = 21450
$ws.Range("M126").Value = -5832.5879
$ws.Range("N126").Value = -26390

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3452.0454
$ws.Range("I132").Value = 2580.2307
$ws.Range("J132").Value = 4711.3335
$ws.Range("K132").Value = 7740.6921
$ws.Range("L132").Value = 14134.0005
$ws.Range("M132").Value = -5210.6921
$ws.Range("N132").Value = -19194.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 797.9143
$ws.Range("I134").Value = 718.8570999999999
$ws.Range("J134").Value = 1114.1428
$ws.Range("K134").Value = 2156.5713
$ws.Range("L134").Value = 3342.4284
$ws.Range("M134").Value = 378.4287000000004
$ws.Range("N134").Value = -8412.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 17796.129
$ws.Range("I136").Value = 1958.75
$ws.Range("J136").Value = 27798.684
$ws.Range("K136").Value = 5876.25
$ws.Range("L136").Value = 83396.052
$ws.Range("M136").Value = -3326.25
$ws.Range("N136").Value = -88496.052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 22416.666
$ws.Range("J137").Value = 23636.363
$ws.Range("L137").Value = 23636.363
$ws.Range("N137").Value = -33836.363

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 11136.556
$ws.Range("I2").Value = 20013.8
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 120082.8
$ws.Range("L2").Value = 240
$ws.Range("M2").Value = -119969.8
$ws.Range("N2").Value = -466

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 737.82355
$ws.Range("I113").Value = 700
$ws.Range("J113").Value = 749.46155
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 2248.38465
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -6588.38465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 967.41174
$ws.Range("I122").Value = 600
$ws.Range("J122").Value = 990.375
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 8913.375
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -13813.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 746.24
$ws.Range("J131").Value = 759.7628999999999
$ws.Range("L131").Value = 2279.2887
$ws.Range("N131").Value = -12359.2887

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1371.4333
$ws.Range("I97").Value = 1455.3478
$ws.Range("J97").Value = 1095.7142
$ws.Range("K97").Value = 1455.3478
$ws.Range("L97").Value = 1095.7142
$ws.Range("M97").Value = -959.3478
$ws.Range("N97").Value = -2087.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3981.3845
$ws.Range("I102").Value = 3467.6365
$ws.Range("K102").Value = 3467.6365
$ws.Range("M102").Value = -1845.6365

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2985.4314
$ws.Range("I126").Value = 2014.3334
$ws.Range("J126").Value = 3848.6296
$ws.Range("K126").Value = 6043.0002
$ws.Range("L126").Value = 11545.8888
$ws.Range("M126").Value = -3573.0002
$ws.Range("N126").Value = -16485.8888

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 24157.076
$ws.Range("I132").Value = 5151.6665
$ws.Range("K132").Value = 15454.9995
$ws.Range("M132").Value = -12924.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 39776
$ws.Range("J135").Value = 39776
$ws.Range("L135").Value = 39776
$ws.Range("N135").Value = -49916

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4915.864
$ws.Range("I7").Value = 4853.0625
$ws.Range("J7").Value = 5083.3335
$ws.Range("K7").Value = 4853.0625
$ws.Range("L7").Value = 5083.3335
$ws.Range("M7").Value = -4741.0625
$ws.Range("N7").Value = -5307.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 1960
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3734.3635
$ws.Range("I40").Value = 3453.9443
$ws.Range("J40").Value = 4996.25
$ws.Range("K40").Value = 3453.9443
$ws.Range("L40").Value = 4996.25
$ws.Range("M40").Value = -3317.9443
$ws.Range("N40").Value = -5268.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2455603
$ws.Range("I122").Value = 2805303.5
$ws.Range("J122").Value = 7700
$ws.Range("K122").Value = 8415910.5
$ws.Range("L122").Value = 23100
$ws.Range("M122").Value = -8413460.5
$ws.Range("N122").Value = -28000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4915.864
$ws.Range("I126").Value = 4853.0625
$ws.Range("J126").Value = 5083.3335
$ws.Range("K126").Value = 14559.1875
$ws.Range("L126").Value = 15250.0005
$ws.Range("M126").Value = -12089.1875
$ws.Range("N126").Value = -20190.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1845.7241
$ws.Range("I136").Value = 1768.5
$ws.Range("K136").Value = 5305.5
$ws.Range("M136").Value = -2755.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 71428880
$ws.Range("I107").Value = 90909310
$ws.Range("K107").Value = 272727930
$ws.Range("M107").Value = -272726010

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1207.3077
$ws.Range("I122").Value = 927.2727
$ws.Range("K122").Value = 2781.8181
$ws.Range("M122").Value = -331.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1946.8334
$ws.Range("I126").Value = 1556.2
$ws.Range("K126").Value = 4668.6
$ws.Range("M126").Value = -2198.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1651.1538
$ws.Range("I132").Value = 1137.7
$ws.Range("K132").Value = 3413.1
$ws.Range("M132").Value = -883.1000000000004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 29496250
$ws.Range("I136").Value = 39703650
$ws.Range("K136").Value = 119110950
$ws.Range("M136").Value = -119108400

Write-Output "Applied all cell updates"